$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 653
$ws.Range("I33").Value = 151.3871
$ws.Range("J33").Value = 4540.5
$ws.Range("K33").Value = 151.3871
$ws.Range("L33").Value = 4540.5
$ws.Range("M33").Value = 77.6129
$ws.Range("N33").Value = -4998.5
$ws.Range("H34").Value = 2709.375
$ws.Range("I34").Value = 953.5714
$ws.Range("J34").Value = 15000
$ws.Range("K34").Value = 953.5714
$ws.Range("L34").Value = 15000
$ws.Range("M34").Value = -750.5714
$ws.Range("N34").Value = -15406
$ws.Range("H36").Value = 2709.375
$ws.Range("I36").Value = 953.5714
$ws.Range("J36").Value = 15000
$ws.Range("K36").Value = 953.5714
$ws.Range("L36").Value = 15000
$ws.Range("M36").Value = -238.5714
$ws.Range("N36").Value = -16430
$ws.Range("H47").Value = 10996
$ws.Range("I47").Value = 5000
$ws.Range("J47").Value = 12495
$ws.Range("K47").Value = 5000
$ws.Range("L47").Value = 12495
$ws.Range("M47").Value = -4028
$ws.Range("N47").Value = -14439
$ws.Range("H51").Value = 7495.227
$ws.Range("J51").Value = 3193.4375
$ws.Range("L51").Value = 3193.4375
$ws.Range("N51").Value = -4161.4375
$ws.Range("H64").Value = 38490.25
$ws.Range("I64").Value = 61462.293
$ws.Range("J64").Value = 2988
$ws.Range("K64").Value = 61462.293
$ws.Range("L64").Value = 2988
$ws.Range("M64").Value = -61214.293
$ws.Range("N64").Value = -3484
$ws.Range("H67").Value = 38490.25
$ws.Range("I67").Value = 61462.293
$ws.Range("J67").Value = 2988
$ws.Range("K67").Value = 61462.293
$ws.Range("L67").Value = 2988
$ws.Range("M67").Value = -60604.293
$ws.Range("N67").Value = -4704
$ws.Range("H74").Value = 3199.923
$ws.Range("I74").Value = 3199.9092
$ws.Range("J74").Value = 3200
$ws.Range("K74").Value = 3199.9092
$ws.Range("L74").Value = 3200
$ws.Range("M74").Value = -2263.9092
$ws.Range("N74").Value = -5072
$ws.Range("H77").Value = 3199.923
$ws.Range("I77").Value = 3199.9092
$ws.Range("J77").Value = 3200
$ws.Range("K77").Value = 15999.546
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = -11319.546
$ws.Range("N77").Value = -25360
$ws.Range("H111").Value = 8609.647000000001
$ws.Range("I111").Value = 11479.333
$ws.Range("J111").Value = 1722.4
$ws.Range("K111").Value = 34437.999
$ws.Range("L111").Value = 5167.200000000001
$ws.Range("M111").Value = -31370.999
$ws.Range("N111").Value = -11301.2
$ws.Range("H135").Value = 2062.875
$ws.Range("I135").Value = 688.3333
$ws.Range("J135").Value = 3830.1428
$ws.Range("K135").Value = 6194.9997
$ws.Range("L135").Value = 34471.2852
$ws.Range("M135").Value = -3659.9997
$ws.Range("N135").Value = -39541.2852

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 935.3555
$ws.Range("I74").Value = 927.4
$ws.Range("J74").Value = 963.2
$ws.Range("K74").Value = 927.4
$ws.Range("L74").Value = 963.2
$ws.Range("M74").Value = -53.39999999999998
$ws.Range("N74").Value = -2711.2
$ws.Range("H77").Value = 935.3555
$ws.Range("I77").Value = 927.4
$ws.Range("J77").Value = 963.2
$ws.Range("K77").Value = 4637
$ws.Range("L77").Value = 4816
$ws.Range("M77").Value = -269
$ws.Range("N77").Value = -13552
$ws.Range("H102").Value = 60687.06
$ws.Range("I102").Value = 144727.14
$ws.Range("J102").Value = 1859
$ws.Range("K102").Value = 144727.14
$ws.Range("L102").Value = 1859
$ws.Range("M102").Value = -143105.14
$ws.Range("N102").Value = -5103
$ws.Range("H122").Value = 2602.125
$ws.Range("I122").Value = 2262.4285
$ws.Range("J122").Value = 4980
$ws.Range("K122").Value = 6787.2855
$ws.Range("L122").Value = 14940
$ws.Range("M122").Value = -4337.2855
$ws.Range("N122").Value = -19840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 465.3913
$ws.Range("I94").Value = 426.94446
$ws.Range("J94").Value = 603.8
$ws.Range("K94").Value = 426.94446
$ws.Range("L94").Value = 603.8
$ws.Range("M94").Value = 24.05554000000001
$ws.Range("N94").Value = -1505.8
$ws.Range("H99").Value = 1724.826
$ws.Range("I99").Value = 1261.8334
$ws.Range("J99").Value = 1888.2354
$ws.Range("K99").Value = 1261.8334
$ws.Range("L99").Value = 1888.2354
$ws.Range("M99").Value = 236.1666
$ws.Range("N99").Value = -4884.2354
$ws.Range("H105").Value = 69448
$ws.Range("I105").Value = 74251.42999999999
$ws.Range("J105").Value = 2200
$ws.Range("K105").Value = 74251.42999999999
$ws.Range("L105").Value = 2200
$ws.Range("M105").Value = -72504.42999999999
$ws.Range("N105").Value = -5694
$ws.Range("H107").Value = 111162936
$ws.Range("I107").Value = 166743170
$ws.Range("J107").Value = 2484.6667
$ws.Range("K107").Value = 166743170
$ws.Range("L107").Value = 2484.6667
$ws.Range("M107").Value = -166741250
$ws.Range("N107").Value = -6324.6667
$ws.Range("H134").Value = 2315.7334
$ws.Range("I134").Value = 2058.32
$ws.Range("J134").Value = 3602.8
$ws.Range("K134").Value = 6174.960000000001
$ws.Range("L134").Value = 10808.4
$ws.Range("M134").Value = -3639.960000000001
$ws.Range("N134").Value = -15878.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 18063.9
$ws.Range("J68").Value = 18063.9
$ws.Range("L68").Value = 18063.9
$ws.Range("N68").Value = -19561.9
$ws.Range("H71").Value = 18063.9
$ws.Range("J71").Value = 18063.9
$ws.Range("L71").Value = 54191.7
$ws.Range("N71").Value = -61679.7
$ws.Range("H99").Value = 7696.5884
$ws.Range("I99").Value = 2047.3077
$ws.Range("J99").Value = 26056.75
$ws.Range("K99").Value = 2047.3077
$ws.Range("L99").Value = 26056.75
$ws.Range("M99").Value = -549.3077000000001
$ws.Range("N99").Value = -29052.75
$ws.Range("H126").Value = 7696.5884
$ws.Range("I126").Value = 2047.3077
$ws.Range("J126").Value = 26056.75
$ws.Range("K126").Value = 6141.9231
$ws.Range("L126").Value = 78170.25
$ws.Range("M126").Value = -3671.9231
$ws.Range("N126").Value = -83110.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 820.41
$ws.Range("J131").Value = 825.6667
$ws.Range("L131").Value = 2477.0001
$ws.Range("N131").Value = -12557.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 38455.5
$ws.Range("J70").Value = 5528.4
$ws.Range("L70").Value = 5528.4
$ws.Range("N70").Value = -6068.4
$ws.Range("H73").Value = 38455.5
$ws.Range("J73").Value = 5528.4
$ws.Range("L73").Value = 5528.4
$ws.Range("N73").Value = -7400.4
$ws.Range("H126").Value = 6538991
$ws.Range("I126").Value = 4473
$ws.Range("J126").Value = 11766605
$ws.Range("K126").Value = 13419
$ws.Range("L126").Value = 35299815
$ws.Range("M126").Value = -10949
$ws.Range("N126").Value = -35304755

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 401.5
$ws.Range("I10").Value = 401.5
$ws.Range("K10").Value = 401.5
$ws.Range("M10").Value = -261.5
$ws.Range("H61").Value = 1850.7858
$ws.Range("I61").Value = 1761.1
$ws.Range("J61").Value = 2075
$ws.Range("K61").Value = 1761.1
$ws.Range("L61").Value = 2075
$ws.Range("M61").Value = -1559.1
$ws.Range("N61").Value = -2479
$ws.Range("H68").Value = 3753.5
$ws.Range("I68").Value = 2115.3333
$ws.Range("J68").Value = 5391.6665
$ws.Range("K68").Value = 2115.3333
$ws.Range("L68").Value = 5391.6665
$ws.Range("M68").Value = -1366.3333
$ws.Range("N68").Value = -6889.6665
$ws.Range("H71").Value = 3753.5
$ws.Range("I71").Value = 2115.3333
$ws.Range("J71").Value = 5391.6665
$ws.Range("K71").Value = 10576.6665
$ws.Range("L71").Value = 26958.3325
$ws.Range("M71").Value = -6832.666499999999
$ws.Range("N71").Value = -34446.3325
$ws.Range("H82").Value = 1841.6666
$ws.Range("I82").Value = 2534.6667
$ws.Range("J82").Value = 1668.4166
$ws.Range("K82").Value = 2534.6667
$ws.Range("L82").Value = 1668.4166
$ws.Range("M82").Value = -2173.6667
$ws.Range("N82").Value = -2390.4166
$ws.Range("H85").Value = 1841.6666
$ws.Range("I85").Value = 2534.6667
$ws.Range("J85").Value = 1668.4166
$ws.Range("K85").Value = 2534.6667
$ws.Range("L85").Value = 1668.4166
$ws.Range("M85").Value = -1286.6667
$ws.Range("N85").Value = -4164.4166
$ws.Range("H93").Value = 1599.9615
$ws.Range("I93").Value = 1627.3889
$ws.Range("J93").Value = 1538.25
$ws.Range("K93").Value = 1627.3889
$ws.Range("L93").Value = 1538.25
$ws.Range("M93").Value = -379.3888999999999
$ws.Range("N93").Value = -4034.25
$ws.Range("H113").Value = 1850.7858
$ws.Range("I113").Value = 1761.1
$ws.Range("J113").Value = 2075
$ws.Range("K113").Value = 1761.1
$ws.Range("L113").Value = 2075
$ws.Range("M113").Value = 408.9000000000001
$ws.Range("N113").Value = -6415

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41248
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126240
$ws.Range("H96").Value = 125001630
$ws.Range("I96").Value = 166668510
$ws.Range("J96").Value = 1000
$ws.Range("K96").Value = 166668510
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = -166667137
$ws.Range("N96").Value = -3746
$ws.Range("H122").Value = 2282.7058
$ws.Range("I122").Value = 1472.5714
$ws.Range("J122").Value = 2849.8
$ws.Range("K122").Value = 4417.7142
$ws.Range("L122").Value = 8549.400000000001
$ws.Range("M122").Value = -1967.7142
$ws.Range("N122").Value = -13449.4
$ws.Range("H132").Value = 1910.2693
$ws.Range("I132").Value = 1873.9762
$ws.Range("J132").Value = 2062.7
$ws.Range("K132").Value = 5621.9286
$ws.Range("L132").Value = 6188.099999999999
$ws.Range("M132").Value = -3091.9286
$ws.Range("N132").Value = -11248.1
$ws.Range("H141").Value = 61151.43
$ws.Range("J141").Value = 61151.43
$ws.Range("L141").Value = 61151.43
$ws.Range("N141").Value = -71511.42999999999
